$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44365
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 25000
$ws.Range("L2").Value = 25000
$ws.Range("M2").Value = 25000
$ws.Range("P2").Value = 1667
$ws.Range("D3").Value = 44792
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 24000
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = 24000
$ws.Range("P3").Value = 1600
$ws.Range("D4").Value = 44754
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 25000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 25000
$ws.Range("P4").Value = 1667
$ws.Range("D5").Value = 44757
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("P5").Value = 1667
$ws.Range("D6").Value = 44771
$ws.Range("J6").Value = 90
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 25000
$ws.Range("P6").Value = 1667
$ws.Range("D7").Value = 44819
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 22000
$ws.Range("P7").Value = 1467
$ws.Range("D8").Value = 44761
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 23000
$ws.Range("M8").Value = 24000
$ws.Range("P8").Value = 1600
$ws.Range("D9").Value = 44400
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 25000
$ws.Range("P9").Value = 1667
$ws.Range("D10").Value = 44806
$ws.Range("J10").Value = 70
$ws.Range("K10").Value = 23000
$ws.Range("L10").Value = 23000
$ws.Range("M10").Value = 23000
$ws.Range("P10").Value = 1533
$ws.Range("D11").Value = 44799
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 23000
$ws.Range("L11").Value = 23000
$ws.Range("M11").Value = 23000
$ws.Range("P11").Value = 1533
$ws.Range("D12").Value = 44750
$ws.Range("D14").Value = 44740
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 25000
$ws.Range("L14").Value = 25000
$ws.Range("M14").Value = 25000
$ws.Range("P14").Value = 1667
$ws.Range("D15").Value = 44764
$ws.Range("K15").Value = 24000
$ws.Range("L15").Value = 24000
$ws.Range("M15").Value = 24000
$ws.Range("P15").Value = 1600
$ws.Range("D18").Value = 44778
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 24000
$ws.Range("L18").Value = 24000
$ws.Range("M18").Value = 24000
$ws.Range("P18").Value = 1600
$ws.Range("D19").Value = 44838
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 22000
$ws.Range("L19").Value = 22000
$ws.Range("M19").Value = 22000
$ws.Range("P19").Value = 1467
$ws.Range("D20").Value = 44817
$ws.Range("J20").Value = 90
$ws.Range("D21").Value = 44789
$ws.Range("K21").Value = 24000
$ws.Range("L21").Value = 24000
$ws.Range("M21").Value = 24000
$ws.Range("P21").Value = 1600
$ws.Range("D22").Value = 44831
$ws.Range("J22").Value = 90
$ws.Range("K22").Value = 25000
$ws.Range("M22").Value = 25000
$ws.Range("P22").Value = 1667
$ws.Range("D23").Value = 44781
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = 24000
$ws.Range("L23").Value = 24000
$ws.Range("M23").Value = 24000
$ws.Range("P23").Value = 1600
$ws.Range("D24").Value = 44810
$ws.Range("J24").Value = 110
$ws.Range("K24").Value = 22000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 22000
$ws.Range("P24").Value = 1467
$ws.Range("D25").Value = 44418
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("M25").Value = 25000
$ws.Range("P25").Value = 1667
$ws.Range("D26").Value = 44407
$ws.Range("J26").Value = 90
$ws.Range("D27").Value = 44827
$ws.Range("J27").Value = 90
$ws.Range("K27").Value = 22000
$ws.Range("L27").Value = 22000
$ws.Range("M27").Value = 22000
$ws.Range("P27").Value = 1467
